# Auto-generated Excel COM-interop script applying the scheduled market-data refresh
# described by the commit diff. For each affected Leve row we: (1) set the
# updated numeric values for currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N), and (2) clear any column whose cell must end up blank
# (no value at all) rather than holding 0 or some stale number.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3777.7778
$ws.Range("I64").Value = 3252.0625
$ws.Range("J64").Value = 4829.2085
$ws.Range("K64").Value = 3252.0625
$ws.Range("L64").Value = 4829.2085
$ws.Range("M64").Value = -3004.0625
$ws.Range("N64").Value = -5325.2085

$ws.Range("H67").Value = 3777.7778
$ws.Range("I67").Value = 3252.0625
$ws.Range("J67").Value = 4829.2085
$ws.Range("K67").Value = 3252.0625
$ws.Range("L67").Value = 4829.2085
$ws.Range("M67").Value = -2394.0625
$ws.Range("N67").Value = -6545.2085

$ws.Range("H98").Value = 1080.8334
$ws.Range("I98").Value = 1080.8334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1080.8334
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 417.1666
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 1080.8334
$ws.Range("I122").Value = 1080.8334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3242.5002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -792.5001999999999
$ws.Range("N122").ClearContents()

$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H129").Value = 1190.6136
$ws.Range("J129").Value = 1302.2894
$ws.Range("L129").Value = 3906.8682
$ws.Range("N129").Value = -13906.8682

$ws.Range("H132").Value = 1301.6428
$ws.Range("I132").Value = 677.25
$ws.Range("J132").Value = 2134.1667
$ws.Range("K132").Value = 2031.75
$ws.Range("L132").Value = 6402.500100000001
$ws.Range("M132").Value = 498.25
$ws.Range("N132").Value = -11462.5001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -16040
$ws.Range("M39").ClearContents()

$ws.Range("H42").Value = 50031
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 50031
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 50031
$ws.Range("N42").Value = -51003
$ws.Range("M42").ClearContents()

$ws.Range("H76").Value = 29933.334
$ws.Range("J76").Value = 29933.334
$ws.Range("L76").Value = 29933.334
$ws.Range("N76").Value = -30609.334

$ws.Range("H79").Value = 29933.334
$ws.Range("J79").Value = 29933.334
$ws.Range("L79").Value = 29933.334
$ws.Range("N79").Value = -32273.334

$ws.Range("H132").Value = 3765.5908
$ws.Range("I132").Value = 2163.923
$ws.Range("J132").Value = 6079.1113
$ws.Range("K132").Value = 6491.768999999999
$ws.Range("L132").Value = 18237.3339
$ws.Range("M132").Value = -3961.768999999999
$ws.Range("N132").Value = -23297.3339


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 397.27777
$ws.Range("I64").Value = 376.6
$ws.Range("J64").Value = 423.125
$ws.Range("K64").Value = 376.6
$ws.Range("L64").Value = 423.125
$ws.Range("M64").Value = -151.6
$ws.Range("N64").Value = -873.125

$ws.Range("H67").Value = 397.27777
$ws.Range("I67").Value = 376.6
$ws.Range("J67").Value = 423.125
$ws.Range("K67").Value = 376.6
$ws.Range("L67").Value = 423.125
$ws.Range("M67").Value = 403.4
$ws.Range("N67").Value = -1983.125


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4460.403
$ws.Range("I31").Value = 867.1539
$ws.Range("J31").Value = 6739.049
$ws.Range("K31").Value = 867.1539
$ws.Range("L31").Value = 6739.049
$ws.Range("M31").Value = -572.1539
$ws.Range("N31").Value = -7329.049

$ws.Range("H34").Value = 4460.403
$ws.Range("I34").Value = 867.1539
$ws.Range("J34").Value = 6739.049
$ws.Range("K34").Value = 867.1539
$ws.Range("L34").Value = 6739.049
$ws.Range("M34").Value = -665.1539
$ws.Range("N34").Value = -7143.049

$ws.Range("H42").Value = 4950
$ws.Range("I42").Value = 4950
$ws.Range("K42").Value = 4950
$ws.Range("M42").Value = -4357

$ws.Range("H50").Value = 49912.5
$ws.Range("J50").Value = 49912.5
$ws.Range("L50").Value = 49912.5
$ws.Range("N50").Value = -51162.5

$ws.Range("H86").Value = 1954.0625
$ws.Range("I86").Value = 2068.9
$ws.Range("J86").Value = 1762.6666
$ws.Range("K86").Value = 2068.9
$ws.Range("L86").Value = 1762.6666
$ws.Range("M86").Value = -945.9000000000001
$ws.Range("N86").Value = -4008.6666

$ws.Range("H89").Value = 1954.0625
$ws.Range("I89").Value = 2068.9
$ws.Range("J89").Value = 1762.6666
$ws.Range("K89").Value = 10344.5
$ws.Range("L89").Value = 8813.333000000001
$ws.Range("M89").Value = -4728.5
$ws.Range("N89").Value = -20045.333


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 338480.3
$ws.Range("I5").Value = 777
$ws.Range("J5").Value = 409162.4
$ws.Range("K5").Value = 2331
$ws.Range("L5").Value = 1227487.2
$ws.Range("M5").Value = -2219
$ws.Range("N5").Value = -1227711.2

$ws.Range("H98").Value = 476.9
$ws.Range("I98").Value = 480.83334
$ws.Range("K98").Value = 1442.50002
$ws.Range("M98").Value = 55.49998000000005

$ws.Range("H122").Value = 486.09375
$ws.Range("I122").Value = 324.14816
$ws.Range("J122").Value = 1360.6
$ws.Range("K122").Value = 2917.33344
$ws.Range("L122").Value = 12245.4
$ws.Range("M122").Value = -467.3334400000003
$ws.Range("N122").Value = -17145.4

$ws.Range("H131").Value = 1511.55
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 1556.9681
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 4670.9043
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -14750.9043

$ws.Range("H135").Value = 338480.3
$ws.Range("I135").Value = 777
$ws.Range("J135").Value = 409162.4
$ws.Range("K135").Value = 6993
$ws.Range("L135").Value = 3682461.6
$ws.Range("M135").Value = -4458
$ws.Range("N135").Value = -3687531.6


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 57554.5
$ws.Range("J74").Value = 57554.5
$ws.Range("L74").Value = 57554.5
$ws.Range("N74").Value = -59426.5

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H77").Value = 57554.5
$ws.Range("J77").Value = 57554.5
$ws.Range("L77").Value = 172663.5
$ws.Range("N77").Value = -182023.5

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H132").Value = 5299.125
$ws.Range("I132").Value = 2126.6667
$ws.Range("K132").Value = 6380.000100000001
$ws.Range("M132").Value = -3850.000100000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1030.0588
$ws.Range("J22").Value = 1102.4286
$ws.Range("L22").Value = 1102.4286
$ws.Range("N22").Value = -1692.4286

$ws.Range("H27").Value = 1030.0588
$ws.Range("J27").Value = 1102.4286
$ws.Range("L27").Value = 1102.4286
$ws.Range("N27").Value = -1316.4286

$ws.Range("H74").Value = 36600
$ws.Range("J74").Value = 36600
$ws.Range("L74").Value = 36600
$ws.Range("N74").Value = -38596

$ws.Range("H77").Value = 36600
$ws.Range("J77").Value = 36600
$ws.Range("L77").Value = 109800
$ws.Range("N77").Value = -119784


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 44994
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 44994
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 44994
$ws.Range("N75").Value = -46866
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 44994
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 44994
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 134982
$ws.Range("N78").Value = -144342
$ws.Range("M78").ClearContents()

$ws.Range("H113").Value = 905.8823
$ws.Range("I113").Value = 1310.1
$ws.Range("J113").Value = 328.42856
$ws.Range("K113").Value = 3930.3
$ws.Range("L113").Value = 985.28568
$ws.Range("M113").Value = -1760.3
$ws.Range("N113").Value = -5325.28568

$ws.Range("H132").Value = 2215.75
$ws.Range("I132").Value = 1320.8948
$ws.Range("J132").Value = 3523.6155
$ws.Range("K132").Value = 3962.6844
$ws.Range("L132").Value = 10570.8465
$ws.Range("M132").Value = -1432.6844
$ws.Range("N132").Value = -15630.8465

